$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '30.169.07'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '  -0.48%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.913.99'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '  -0.75%  '
$ws.Range('E4').Value = '  -0.27%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.7398'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -1.23%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '244.17'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +0.14%  '
$ws.Range('E7').Value = '  -0.31%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3127'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  -0.77%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '26.88'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  -1.74%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.06998'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +0.57%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.7796'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +1.30%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.08003'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  +0.02%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.911.07'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  -0.82%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.303'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -0.16%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '92.37'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -0.71%  '
$ws.Range('E16').Value = '  +1.06%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '30.169.23'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  -0.42%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '5.929'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  +3.55%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '242.43'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -3.19%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.000007850'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -0.32%  '
$ws.Range('B21').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C21').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '2.175.36'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -0.21%  '
$ws.Range('B22').Value = 'Dai'
$ws.Range('C22').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.000'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -0.20%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.001'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -0.17%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '7.215'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +8.98%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.448'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +0.44%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '168.18'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +1.57%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '19.09'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +1.15%  '
$ws.Range('E28').Value = '  -2.74%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.067'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -5.38%  '
$ws.Range('E30').Value = '  -0.93%  '
$ws.Range('E31').Value = '  +2.41%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.357'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -0.06%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.111'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +0.72%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.05194'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +1.75%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.302'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +2.26%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7536'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +1.62%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.726'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -1.82%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01945'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +0.14%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.806'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +0.33%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '6.391'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +0.88%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '75.33'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -1.83%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.4522'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +2.38%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.968'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +1.13%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '7.851'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +5.88%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.002'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -0.11%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.8394'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +1.05%  '
$ws.Range('B47').Value = 'Quant'
$ws.Range('C47').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '101.96'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +1.54%  '
$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '9.965'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +3.09%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '37.24'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +0.31%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.065.96'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -0.83%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.1201'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +3.42%  '
